# Remove "null" strings from test data.
# The source workbook had several cells whose literal value was the
# string "null" (a placeholder that leaked into the exported test
# fixture). Clearing those cells drops the now-unused "null" entry
# from the shared-strings table on save, which is what the upstream
# commit's diff shows (uniqueCount 45 -> 44, count 75 -> 69, and every
# shared-string index above the old "null" slot shifting down by one).

$wb = $excel.ActiveWorkbook

$wsBasic = $wb.Worksheets.Item("basic info")
$wsConditions = $wb.Worksheets.Item("conditions")
$wsMoreConditions = $wb.Worksheets.Item("more conditions")

# Drop every cell that literally holds "null".
$wsConditions.Range("E2").ClearContents() | Out-Null

$wsMoreConditions.Range("E2").ClearContents() | Out-Null
$wsMoreConditions.Range("B4").ClearContents() | Out-Null
$wsMoreConditions.Range("B5").ClearContents() | Out-Null
$wsMoreConditions.Range("B6").ClearContents() | Out-Null
$wsMoreConditions.Range("E6").ClearContents() | Out-Null

# Re-create the selection/active-sheet state left behind by the edit:
# "basic info" keeps its prior selection, "conditions" ends up selected
# on D12, and "more conditions" becomes the active tab, selected on C6.
$wsBasic.Activate() | Out-Null
$wsBasic.Range("C5").Select() | Out-Null

$wsConditions.Activate() | Out-Null
$wsConditions.Range("D12").Select() | Out-Null

$wsMoreConditions.Activate() | Out-Null
$wsMoreConditions.Range("C6").Select() | Out-Null
